# Inserts a new first data row (DTLANC=46031, TOTAL_BONUS=32) above the
# existing data in the "BONUS" sheet, shifting all prior data rows down by
# one row (A2:B169 -> A3:B170), and extending the used range accordingly
# (A1:B169 -> A1:B170).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing data block (rows 2-169, columns A-B) with raw values
# (Value2 avoids Excel re-interpreting date serials and auto-assigning a
# locale date format to newly created cells).
$data = $ws.Range("A2:B169").Value2()

# Shift the whole block down by one row.
$ws.Range("A3:B170").Value2 = $data

# Write the new row of data into the vacated row 2.
$ws.Range("A2").Value2 = 46031
$ws.Range("B2").Value2 = 32

# The bulk Value2 write leaves the brand-new row 170 with the default
# (unstyled) cell format. Re-apply the date-serial number format used by
# the rest of column A by copying formatting from the row above it.
$ws.Range("A169").Copy()
$ws.Range("A170").PasteSpecial(-4122)
$excel.CutCopyMode = 0
